$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style from A269 (existing last date cell) to the new date column cells A270:A301
$styleSrc = $ws.Cells.Item(269, 1)
$styleSrc.Copy()
$dstStyleRange = $ws.Range("A270:A301")
$dstStyleRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate data for rows 270-301
$ws.Cells.Item(270, 1).Value = 44344
$ws.Cells.Item(270, 2).Value = 4
$ws.Cells.Item(270, 3).Value = 23
$ws.Cells.Item(270, 4).Value = 31.78156392930675
$ws.Cells.Item(271, 1).Value = 44345
$ws.Cells.Item(271, 2).Value = 5
$ws.Cells.Item(271, 3).Value = 22
$ws.Cells.Item(271, 4).Value = 30.39975680194558
$ws.Cells.Item(272, 1).Value = 44346
$ws.Cells.Item(272, 2).Value = 2
$ws.Cells.Item(272, 3).Value = 24
$ws.Cells.Item(272, 4).Value = 33.16337105666791
$ws.Cells.Item(273, 1).Value = 44347
$ws.Cells.Item(273, 2).Value = 1
$ws.Cells.Item(273, 3).Value = 18
$ws.Cells.Item(273, 4).Value = 24.87252829250093
$ws.Cells.Item(274, 1).Value = 44348
$ws.Cells.Item(274, 2).Value = 0
$ws.Cells.Item(274, 3).Value = 15
$ws.Cells.Item(274, 4).Value = 20.72710691041744
$ws.Cells.Item(275, 1).Value = 44349
$ws.Cells.Item(275, 2).Value = 0
$ws.Cells.Item(275, 3).Value = 13
$ws.Cells.Item(275, 4).Value = 17.96349265569512
$ws.Cells.Item(276, 1).Value = 44350
$ws.Cells.Item(276, 2).Value = 1
$ws.Cells.Item(276, 3).Value = 13
$ws.Cells.Item(276, 4).Value = 17.96349265569512
$ws.Cells.Item(277, 1).Value = 44351
$ws.Cells.Item(277, 2).Value = 4
$ws.Cells.Item(277, 3).Value = 13
$ws.Cells.Item(277, 4).Value = 17.96349265569512
$ws.Cells.Item(278, 1).Value = 44352
$ws.Cells.Item(278, 2).Value = 0
$ws.Cells.Item(278, 3).Value = 8
$ws.Cells.Item(278, 4).Value = 11.0544570188893
$ws.Cells.Item(279, 1).Value = 44353
$ws.Cells.Item(279, 2).Value = 2
$ws.Cells.Item(279, 3).Value = 8
$ws.Cells.Item(279, 4).Value = 11.0544570188893
$ws.Cells.Item(280, 1).Value = 44354
$ws.Cells.Item(280, 2).Value = 3
$ws.Cells.Item(280, 3).Value = 10
$ws.Cells.Item(280, 4).Value = 13.81807127361163
$ws.Cells.Item(281, 1).Value = 44355
$ws.Cells.Item(281, 2).Value = 1
$ws.Cells.Item(281, 3).Value = 11
$ws.Cells.Item(281, 4).Value = 15.19987840097279
$ws.Cells.Item(282, 1).Value = 44356
$ws.Cells.Item(282, 2).Value = 0
$ws.Cells.Item(282, 3).Value = 11
$ws.Cells.Item(282, 4).Value = 15.19987840097279
$ws.Cells.Item(283, 1).Value = 44357
$ws.Cells.Item(283, 2).Value = 2
$ws.Cells.Item(283, 3).Value = 12
$ws.Cells.Item(283, 4).Value = 16.58168552833395
$ws.Cells.Item(284, 1).Value = 44358
$ws.Cells.Item(284, 2).Value = 0
$ws.Cells.Item(284, 3).Value = 8
$ws.Cells.Item(284, 4).Value = 11.0544570188893
$ws.Cells.Item(285, 1).Value = 44359
$ws.Cells.Item(285, 2).Value = 3
$ws.Cells.Item(285, 3).Value = 11
$ws.Cells.Item(285, 4).Value = 15.19987840097279
$ws.Cells.Item(286, 1).Value = 44360
$ws.Cells.Item(286, 2).Value = 0
$ws.Cells.Item(286, 3).Value = 9
$ws.Cells.Item(286, 4).Value = 12.43626414625047
$ws.Cells.Item(287, 1).Value = 44361
$ws.Cells.Item(287, 2).Value = 1
$ws.Cells.Item(287, 3).Value = 7
$ws.Cells.Item(287, 4).Value = 9.672649891528142
$ws.Cells.Item(288, 1).Value = 44362
$ws.Cells.Item(288, 2).Value = 0
$ws.Cells.Item(288, 3).Value = 6
$ws.Cells.Item(288, 4).Value = 8.290842764166976
$ws.Cells.Item(289, 1).Value = 44363
$ws.Cells.Item(289, 2).Value = 1
$ws.Cells.Item(289, 3).Value = 7
$ws.Cells.Item(289, 4).Value = 9.672649891528142
$ws.Cells.Item(290, 1).Value = 44364
$ws.Cells.Item(290, 2).Value = 2
$ws.Cells.Item(290, 3).Value = 7
$ws.Cells.Item(290, 4).Value = 9.672649891528142
$ws.Cells.Item(291, 1).Value = 44365
$ws.Cells.Item(291, 2).Value = 0
$ws.Cells.Item(291, 3).Value = 7
$ws.Cells.Item(291, 4).Value = 9.672649891528142
$ws.Cells.Item(292, 1).Value = 44366
$ws.Cells.Item(292, 2).Value = 3
$ws.Cells.Item(292, 3).Value = 7
$ws.Cells.Item(292, 4).Value = 9.672649891528142
$ws.Cells.Item(293, 1).Value = 44367
$ws.Cells.Item(293, 2).Value = 0
$ws.Cells.Item(293, 3).Value = 7
$ws.Cells.Item(293, 4).Value = 9.672649891528142
$ws.Cells.Item(294, 1).Value = 44368
$ws.Cells.Item(294, 2).Value = 2
$ws.Cells.Item(294, 3).Value = 8
$ws.Cells.Item(294, 4).Value = 11.0544570188893
$ws.Cells.Item(295, 1).Value = 44369
$ws.Cells.Item(295, 2).Value = 0
$ws.Cells.Item(295, 3).Value = 8
$ws.Cells.Item(295, 4).Value = 11.0544570188893
$ws.Cells.Item(296, 1).Value = 44370
$ws.Cells.Item(296, 2).Value = 1
$ws.Cells.Item(296, 3).Value = 8
$ws.Cells.Item(296, 4).Value = 11.0544570188893
$ws.Cells.Item(297, 1).Value = 44371
$ws.Cells.Item(297, 2).Value = 0
$ws.Cells.Item(297, 3).Value = 6
$ws.Cells.Item(297, 4).Value = 8.290842764166976
$ws.Cells.Item(298, 1).Value = 44372
$ws.Cells.Item(298, 2).Value = 0
$ws.Cells.Item(298, 3).Value = 6
$ws.Cells.Item(298, 4).Value = 8.290842764166976
$ws.Cells.Item(299, 1).Value = 44373
$ws.Cells.Item(299, 2).Value = 0
$ws.Cells.Item(299, 3).Value = 3
$ws.Cells.Item(299, 4).Value = 4.145421382083488
$ws.Cells.Item(300, 1).Value = 44374
$ws.Cells.Item(300, 2).Value = 0
$ws.Cells.Item(300, 3).Value = 3
$ws.Cells.Item(300, 4).Value = 4.145421382083488
$ws.Cells.Item(301, 1).Value = 44375
$ws.Cells.Item(301, 2).Value = 0
$ws.Cells.Item(301, 3).Value = 1
$ws.Cells.Item(301, 4).Value = 1.381807127361163

Write-Host "Added rows 270-301"
